# WBS uppdaterad, Aktivitetslista och Riskanalys skapade och färdigställda
#
# Adds a "Beroenden" (Dependencies) column (G) to the WBS sheet, fills in the
# dependency references for each task row, removes the stray estimate
# figures that had been entered on row 12 (the "1.2.2 Säkerhet" summary row
# shouldn't carry its own Resurs/Tidsåtgång/Total numbers - those live on its
# sub-rows), and updates the running totals accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Beroenden" column header -----------------------------------------
$ws.Range("G1").Value = "Beroenden"

# --- Dependency values for each task row ------------------------------------
$ws.Range("G4").Value  = "1.1.5, 3.1"
$ws.Range("G5").Value  = "1.1.5, 3.1"
$ws.Range("G6").Value  = "1.1.5, 3.1"
$ws.Range("G7").Value  = "1.1.5, 3.1"
$ws.Range("G8").Value  = "-"
$ws.Range("G9").Value  = "1.3.1, 1.3.2, 1.3.3"
$ws.Range("G11").Value = "1.1.1, 1.1.2, 1.1.3, 1.1.4"
$ws.Range("G13").Value = "1.2.1"
$ws.Range("G14").Value = "1.2.1"
$ws.Range("G15").Value = "1.2.1"
$ws.Range("G16").Value = "1.2.1"
$ws.Range("G17").Value = "1.2.1"
$ws.Range("G19").Value = "-"
$ws.Range("G20").Value = "-"
$ws.Range("G21").Value = "-"
$ws.Range("G23").Value = "1.1.6"
$ws.Range("G26").Value = 1.4
$ws.Range("G27").Value = "-"
$ws.Range("G28").Value = "2.1.1"
$ws.Range("G29").Value = "2.1.1"
$ws.Range("G30").Value = 2
$ws.Range("G31").Value = "2.1.1"
$ws.Range("G32").Value = "2.1.1"
$ws.Range("G33").Value = "2.1.1"
$ws.Range("G35").Value = "1-"
$ws.Range("G36").Value = "-"
$ws.Range("G37").Value = 1.4

# --- Remove the leftover estimate numbers on the 1.2.2 summary row ---------
$ws.Range("D12:F12").ClearContents()

# --- Formatting ---------------------------------------------------------
# Column A already uses a "left" style (s=1); the new dependency column
# uses a "right" aligned style, matching Excel's auto-generated cellXfs.
$ws.Range("G1:G37").HorizontalAlignment = -4152
$ws.Columns.Item(6).ColumnWidth = 25
$ws.Columns.Item(7).ColumnWidth = 20.33203125

# --- View state: scroll down a bit and park the active cell/selection -------
$ws.Application.Goto($ws.Range("A11"), $true)
$ws.Range("H12").Select()
